$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for the 380 kV case (rows 2-25, columns B-F, I-N).
# Column order matches the sheet layout: B,C,D,E,F,I,J,K,L,M,N (G and H unaffected).
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.04875610185102, 1.056938332385505, 1.05639297797345, 1.067956458803336, 1.048502342011862, 1.053798165644134, 1.059674671496669, 1.059130812873581, 1.070662931607358, 1.055294679911688)
    3 = @(1.02, 1.049659615628725, 1.057675328773298, 1.05719811417175, 1.068846106531008, 1.048764196934572, 1.054350910187309, 1.060225600802422, 1.059749600776865, 1.071368281079545, 1.055848209415514)
    4 = @(1.02, 1.050244775106865, 1.058152638480357, 1.057719929134223, 1.069422709596494, 1.048932560864472, 1.054708438807389, 1.060581842744047, 1.060150179424004, 1.071824989814458, 1.056206245767244)
    5 = @(1.02, 1.050490900668912, 1.058353399354572, 1.057939499160407, 1.069665337236336, 1.04900308328974, 1.054858710611272, 1.060731546720789, 1.060318624868923, 1.072017060673677, 1.056356730974326)
    6 = @(1.02, 1.050532233497152, 1.058387113799534, 1.057976377557721, 1.069706088554234, 1.049014909178518, 1.054883939926932, 1.060756679115835, 1.060346910036996, 1.072049314309282, 1.056381996118507)
    7 = @(1.02, 1.050248063361052, 1.058155320666528, 1.057722862259675, 1.069425950721973, 1.048933504202332, 1.054710446878569, 1.060583843333711, 1.060152430036129, 1.071827555999419, 1.056208256690116)
    8 = @(1.02, 1.04906133892668, 1.057187315074991, 1.056664903042683, 1.068256923461836, 1.048591058994543, 1.053984995326663, 1.059860910879378, 1.059339896799512, 1.07090124469521, 1.055481774913796)
    9 = @(1.02, 1.046974263232263, 1.055484880686516, 1.054807132716529, 1.066204233283889, 1.047979435321469, 1.052705680359822, 1.058585179472946, 1.057909555544189, 1.069271329924838, 1.054200643172962)
    10 = @(1.02, 1.045585703723154, 1.054352246896, 1.053573077247953, 1.064840771374093, 1.047566221481049, 1.051852214241733, 1.05773353150653, 1.056957042271201, 1.068186394043908, 1.053345965035092)
    11 = @(1.02, 1.044985128245637, 1.053862374450024, 1.053039795003069, 1.064251583807707, 1.047386008078181, 1.051482528212942, 1.057364498251745, 1.056544857743757, 1.067717019571943, 1.052975754009736)
    12 = @(1.02, 1.044762151204056, 1.053680500537419, 1.05284187288247, 1.064032915347019, 1.047318875816605, 1.051345192006257, 1.057227384380065, 1.056391794376677, 1.067542735992056, 1.052838222769884)
    13 = @(1.02, 1.044809975835366, 1.053719509168376, 1.05288432046318, 1.064079812209631, 1.04733328463623, 1.05137465190156, 1.057256797493715, 1.056424625143811, 1.067580117538103, 1.052867724501618)
    14 = @(1.02, 1.044966694776454, 1.053847338921491, 1.053023431371546, 1.064233504881337, 1.047380462837222, 1.051471176327267, 1.057353165154667, 1.056532204648323, 1.067702611943556, 1.052964386003082)
    15 = @(1.02, 1.045063268232351, 1.053926110563754, 1.053109163794115, 1.06432822418063, 1.047409505332935, 1.051530645812645, 1.057412535386273, 1.056598493311858, 1.067778093179648, 1.053023939941949)
    16 = @(1.02, 1.045625576374307, 1.054384770164621, 1.053608492155384, 1.064879899247912, 1.047578154550461, 1.051876746450877, 1.057758017565691, 1.056984403186007, 1.068217553638855, 1.053370532082787)
    17 = @(1.02, 1.045978480358041, 1.054672627713765, 1.053921995760788, 1.065226273021703, 1.04768359911149, 1.052093812332209, 1.057974659546516, 1.057226544751802, 1.068493326318993, 1.053587906222568)
    18 = @(1.02, 1.046184388800887, 1.054840584605309, 1.054104960426462, 1.065428422695759, 1.047744978715699, 1.052220410612372, 1.05810099759077, 1.0573678068329, 1.068654219314636, 1.05371468428681)
    19 = @(1.02, 1.04625460933328, 1.054897862753336, 1.054167364124662, 1.065497370086168, 1.04776588644906, 1.052263575187442, 1.058144071227953, 1.057415977754667, 1.068709086317969, 1.053757910160527)
    20 = @(1.02, 1.045940610314184, 1.054641737679284, 1.053888349092762, 1.065189098413032, 1.047672298772804, 1.052070524510193, 1.057951418547271, 1.057200562658196, 1.068463734428876, 1.053564585329173)
    21 = @(1.02, 1.044920542088938, 1.05380969384413, 1.052982462173459, 1.064188241181634, 1.047366575354546, 1.051442752787475, 1.057324788333415, 1.056500524042503, 1.067666538643854, 1.052935922098603)
    22 = @(1.02, 1.044279783766466, 1.053287056248101, 1.052413836947801, 1.063560016548426, 1.04717323872951, 1.051047942141761, 1.056930578944283, 1.056060615875532, 1.0671656751989, 1.052540550776486)
    23 = @(1.02, 1.044619404665774, 1.053564068276096, 1.05271518612731, 1.063892949840767, 1.04727583570211, 1.051257248378056, 1.057139577457594, 1.056293796928033, 1.067431157300825, 1.052750154251644)
    24 = @(1.02, 1.045957721954583, 1.05465569538815, 1.053903552254696, 1.065205895661191, 1.047677405293957, 1.052081047313525, 1.057961920233531, 1.057212302770796, 1.068477105608494, 1.053575123076093)
    25 = @(1.02, 1.047513330264959, 1.055924598279544, 1.055286632082135, 1.066734029255496, 1.048138521034439, 1.053036522763986, 1.058915196373465, 1.058279153916614, 1.069692413101446, 1.054531955411288)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowVals[$i]
    }
}
